$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.595.16'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.68%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.821.70'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.41%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.009'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.21%  '

$ws.Range('E5').Value = '  +0.18%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '305.29'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.72%  '

$ws.Range('E7').Value = '  +2.10%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3592'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07126'
$ws.Range('D9').ClearFormats()

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8968'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.20%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07773'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.63%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.30'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.09%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.828.24'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.73%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.244'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.82%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.322'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.17%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.20'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.37%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.010'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.26%  '

$ws.Range('E18').Value = '  -0.73%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.008'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.16%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.642.67'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.78%  '

$ws.Range('E21').Value = '  -1.14%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.010'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.38%  '

$ws.Range('E23').Value = '  -0.40%  '

$ws.Range('E24').Value = '  -3.57%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.06'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.47%  '

$ws.Range('E26').Value = '  -0.21%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.964'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -4.19%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '113.53'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.37%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.792'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.67%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08791'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.53%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.132'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.40%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7276'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.05%  '

$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.723'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.21%  '

$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.423'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.67%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.125'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.17%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.074'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.14%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01919'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.35%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.917'
$ws.Range('D38').ClearFormats()

$ws.Range('E39').Value = '  -0.33%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.843'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.79%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5031'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.78%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1490'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.23%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.945'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.88%  '

$ws.Range('E44').Value = '  +0.20%  '

$ws.Range('E45').Value = '  -0.82%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.909'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.58%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '98.06'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.14%  '

$ws.Range('E48').Value = '  -2.27%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05981'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.07%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '63.65'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.91%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.71'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.73%  '
